# GaN_2022_ActivityGuide_Orion (Chinese Traditional) — consolidate the
# multi-run "shorttext" / paragraph / links text into single runs that
# use the document's existing GaNStyle / GaNParagraph / GaNLinks
# character styles, updating the constellation/date text along the way.
#
# Pattern used for every paragraph being consolidated:
#   1. Find a unique anchor string to locate the paragraph.
#   2. Expand the found range to the whole paragraph (wdParagraph) and
#      pull the end back one character so the paragraph mark / section
#      break mark is excluded.
#   3. Delete the old (multi-run) content — this collapses the range.
#   4. Insert the new consolidated text into the now-empty range; since
#      the range is collapsed, the new text does NOT inherit the old
#      runs' direct formatting (it only gets the paragraph's default
#      run formatting).
#   5. Apply the target character style to the freshly inserted range.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) The four "獵戶座 / 英仙座" short red headline paragraphs.
#    They all collapse to the same new sentence with style "GaNStyle".
# ---------------------------------------------------------------------
$orionText = "獵戶座： 2022年1月16-25日、2月14-23日、3月14-24日."

$rng = $d.Content
while ($rng.Find.Execute("英仙座：")) {
    $rng.Expand(4)
    $rng.MoveEnd(1, -1)
    $rng.Delete()
    $rng.InsertAfter($orionText)
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# ---------------------------------------------------------------------
# 2) The "你现在参加的是全球公益科普活动 ..." intro paragraph —
#    collapses to a single run with style "GaNParagraph".
# ---------------------------------------------------------------------
$introText = "你现在参加的是全球公益科普活动 Globe at Night （全球观星活动），这是一个以观察和记录夜空的可见恒星数来测量你所在地光污染情况的活动。通过定位和观测夜空中的獵戶座，并将所肉眼观察到的星等情况与所给出的星等图表作对比，我们可以知道自己社区中的人造光是如何导致光污染的。你所提供数据将和来自全世界的数据一起帮助建立一张全球光污染地图。"

$rng = $d.Content
if ($rng.Find.Execute("你现在参加的是")) {
    $rng.Expand(4)
    $rng.MoveEnd(1, -1)
    $rng.Delete()
    $rng.InsertAfter($introText)
    $rng.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------
# 3) The "本文檔中的圖表由 Jenik Hollan, CzechGlobe (...)." credit line
#    — collapses to a single run with style "GaNLinks" (also bumps the
#    map year from 2019 to 2022).
# ---------------------------------------------------------------------
$linksText = "本文檔中的圖表由 Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
if ($rng.Find.Execute("本文檔中的圖表由")) {
    $rng.Expand(4)
    $rng.MoveEnd(1, -1)
    $rng.Delete()
    $rng.InsertAfter($linksText)
    $rng.Style = "GaNLinks"
}

Write-Output "done"
